# Update scripts with new TPM values.
#
# The "Sending cluster" for the first data row changes from FAPs -> ECs,
# and a new second data row for FAPs (previously Resolving-Mac) replaces
# the old Resolving-Mac row, with all of its NATMI-derived metrics
# recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster ECs -> Pnoc -> Oprl1 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1534156666666667
$ws.Range("H2").Value = 0.460247
$ws.Range("I2").Value = 0.4248984253037086
$ws.Range("J2").Value = 0.4248984253037086
$ws.Range("Q2").Value = 0.1113364596434444
$ws.Range("R2").Value = 1.002028136791
$ws.Range("S2").Value = 0.4248984253037086
$ws.Range("T2").Value = 0.4248984253037086

# Row 3: Sending cluster FAPs (was Resolving-Mac) -> Pnoc -> Oprl1 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.2076486666666667
$ws.Range("H3").Value = 0.622946
$ws.Range("I3").Value = 0.5751015746962914
$ws.Range("J3").Value = 0.5751015746962914
$ws.Range("Q3").Value = 0.1506943058597778
$ws.Range("R3").Value = 1.356248752738
$ws.Range("S3").Value = 0.5751015746962914
$ws.Range("T3").Value = 0.5751015746962914
